# Update imputed values in result_data_RandomForest sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 16.08500000000002
$ws.Range("C9").Value = -10.47480000000001
$ws.Range("E12").Value = 18.02740000000003
$ws.Range("C13").Value = -12.62199999999999
$ws.Range("E14").Value = 16.54610000000001
$ws.Range("C16").Value = -13.0579
$ws.Range("C18").Value = -12.5426
$ws.Range("E19").Value = 16.6063
$ws.Range("C20").Value = -11.64980000000001
$ws.Range("C26").Value = -13.02500000000001
$ws.Range("E26").Value = 15.9293
$ws.Range("C27").Value = -12.22889999999999
$ws.Range("E27").Value = 16.64379999999999
$ws.Range("C29").Value = -11.4927
$ws.Range("E29").Value = 17.08920000000001
$ws.Range("C35").Value = -12.73920000000001
$ws.Range("C36").Value = -12.85920000000001
$ws.Range("E37").Value = 16.66760000000001
$ws.Range("E38").Value = 16.409
$ws.Range("C45").Value = -13.62419999999999
$ws.Range("E47").Value = 16.5685
$ws.Range("E51").Value = 17.1259
$ws.Range("E52").Value = 17.1937
$ws.Range("C55").Value = -13.7679
$ws.Range("E55").Value = 16.41290000000001
$ws.Range("C57").Value = -13.4868
$ws.Range("C69").Value = -11.7344
$ws.Range("E69").Value = 17.05590000000002
$ws.Range("E70").Value = 17.91610000000003
$ws.Range("C76").Value = -12.1712
$ws.Range("E76").Value = 16.6507
$ws.Range("C78").Value = -11.84860000000001
$ws.Range("E81").Value = 16.3624
$ws.Range("C82").Value = -12.0795
$ws.Range("C83").Value = -14.059
$ws.Range("E83").Value = 16.7029
$ws.Range("C93").Value = -10.7465
$ws.Range("E94").Value = 18.24060000000002
$ws.Range("C97").Value = -12.8712
$ws.Range("E100").Value = 16.5713
$ws.Range("E102").Value = 16.8129
